$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("D3").Value = 10.65
$ws.Range("F3").Value = 10.05

$ws.Range("C4").Value = 9.35
$ws.Range("E4").Value = 9.949999999999999
$ws.Range("F4").Value = 10.02

$ws.Range("D5").Value = 10.05
$ws.Range("F5").Value = 10.05
$ws.Range("G5").Value = 9.32

$ws.Range("C6").Value = 9.949999999999999
$ws.Range("D6").Value = 9.98
$ws.Range("E6").Value = 9.949999999999999
$ws.Range("G6").Value = 10.11
$ws.Range("J6").Value = 7.43

$ws.Range("E7").Value = 10.68
$ws.Range("F7").Value = 9.890000000000001

$ws.Range("J9").Value = 12.67

$ws.Range("F10").Value = 12.57
$ws.Range("I10").Value = 7.33
